$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4; this shifts the existing rows 4-38
# (and their formatting) down to rows 5-39, extending the data range
# from A1:T38 to A1:T39.
$ws.Rows("4:4").Insert()

# Populate the newly inserted row 4 with the new weekly price entry.
$newRow = @(
    4,
    "Feria Lagunitas de Puerto Montt",
    "Los Lagos",
    44530,
    10,
    "Fruta",
    100103,
    "Frutos de hueso (carozo)",
    100103001,
    "Cereza",
    "Santina",
    "Primera",
    600,
    20000,
    21000,
    20500,
    "`$/caja 15 kilos",
    "Provincia de Curicó",
    1367,
    15
)

for ($i = 0; $i -lt $newRow.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $newRow[$i]
}
